$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19 (shifts existing rows 19-44 down to 20-45)
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly record
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C19").Value = "Arica y Parinacota"
$ws.Range("D19").Value2 = 44883
$ws.Range("E19").Value = 15
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100104
$ws.Range("H19").Value = "Frutos de pepita"
$ws.Range("I19").Value = 100104005
$ws.Range("J19").Value = "Pera"
$ws.Range("K19").Value = "Packham's Triumph"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 300
$ws.Range("N19").Value = 24000
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 24500
$ws.Range("Q19").Value = "$/bandeja 18 kilos granel"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 1361
$ws.Range("T19").Value = 18
